$d = $word.ActiveDocument

# Map of old text -> new text for the four fields being updated.
$replacements = @{
    "Test HOspital"                 = "Frontline Hospital 1"
    "s.sudippaudel@gmail.com"       = "polo@yopmail.com"
    "February 5, 2024, 10:00 AM"    = "March 30, 2024, 8:22 PM"
}

# The hospital name, e-mail and report date/time all live inside content
# controls (structured document tags) in this template, so Document.Content.Find
# does not reach them directly -- walk the ContentControls collection and
# match on the current text instead.
$ccs = $d.ContentControls
for ($i = 1; $i -le $ccs.Count; $i++) {
    $cc = $ccs.Item($i)
    $current = $cc.Range.Text
    if ($replacements.ContainsKey($current)) {
        $cc.Range.Text = $replacements[$current]
    }
}

# The doctor's name is plain body text (not inside a content control), so a
# normal Find/Replace on the document body works for it.
$null = $d.Content.Find.Execute("Rhythm  Sapkota", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Hernam Lopchan", 2)
